$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "BB"
$ws.Range("B3").Value = "Bangalore"

$ws.Range("D3").Select()
